$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()
$ws.Range("I19").Select()
$win = $excel.ActiveWindow
$win.Split = $true
$ws.Range("D47").Select()
